$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8 ("life cycle" / $T$ / 40 / standard calibration).
# This shifts all rows from 8 downward by one, matching the diff (rows 8-21 -> 9-22).
$ws.Rows(8).Insert()

# Copy formatting (style, borders, etc.) from the row that is now below (old row 8, now row 9)
# into the freshly inserted row 8 so the new row looks consistent with the rest of the table.
$ws.Range("A9:D9").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)

# Populate the new row with the population growth parameter.
$ws.Range("A8").Value = "life cycle"
$ws.Range("B8").Value = "`$n`$"
$ws.Range("C8").Value = 0.005
$ws.Range("D8").Value = "U.S. Census"

Write-Host "Inserted row 8: life cycle / `$n`$ / 0.005 / U.S. Census"
